# Adding the changes we made on may 9th
#
# Inserts 7 new data rows before the existing data (pushing the old
# data down from rows 2-21 to rows 9-28), and appends 3 new data rows
# at the end (rows 29-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 7 blank rows above the current row 2 (the first
# data row), shifting the existing 20 data rows down to rows 9-28. ---
$ws.Range("A2:C8").Insert()
# Row insertion copies formatting from the row above (the bold header
# row); clear that so the new cells stay unstyled like the rest of the
# data rows.
$ws.Range("A2:C8").ClearFormats()

# --- Step 2: populate the 7 newly inserted rows (2-8) ---
$topData = @(
    @(-0.0310014113783836, 0.0154243474826216, 0.0794124826788902),
    @(-0.005192354787141, 0.0403171069920063, -0.0171042270958423),
    @(0.0421497002243995, -0.0251981914043426, 0.0042760567739605),
    @(0.0103847095742821, -0.0126754539087414, -0.00335975876078),
    @(0.0180205255746841, -0.0375682115554809, 0.0152716310694813),
    @(-0.0210748501121997, -0.0261144898831844, -0.0280998013913631),
    @(-0.0584903471171855, 0.0059559359215199, 0.0522289797663688)
)

$r = 2
foreach ($row in $topData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Step 3: append 3 new data rows at the bottom (rows 29-31) ---
$bottomData = @(
    @(-0.0070249503478407, -0.0122173046693205, 0.0006108652451075),
    @(-0.0004581489483825, 0.0073303831741213, 0.0157297793775796),
    @(-0.0102319931611418, -0.020616702735424, 0.0047342055477201)
)

$r = 29
foreach ($row in $bottomData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
